$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update designator list & quantity for the 2POS header part (row 3)
$ws.Range("C3").Value = "'AUXM, DRVM, LEDB, PROP, SPKR, STRM, SW, TRB"
$ws.Range("F3").Value = 8

# Remove the obsolete SW400204-1 / DRVM, SPKR row entirely (row 15);
# remaining rows shift up by one
$ws.Rows(15).Delete()

# Row 16 (was 17): Ferrite comment -> actual part number
$ws.Range("A16").Value = "'MPZ1608S221ATA00"

# Row 17 (was 18): swap in the -R7 tape/reel variant part number
$ws.Range("A17").Value = "'AD623ARZ-R7"

# Row 31 (was 32): LibRef corrected from CS326 to INA219BIDR
$ws.Range("E31").Value = "'INA219BIDR"
